$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1868512110726644
$ws.Range("C2").Value = 0.5536332179930796
$ws.Range("J2").Value = 0.01384083044982699
$ws.Range("P2").Value = 0.1314878892733564
$ws.Range("S2").Value = 0.1141868512110727
$ws.Range("B3").Value = 0.04191616766467066
$ws.Range("C3").Value = 0.02395209580838323
$ws.Range("J3").Value = 0.03592814371257485
$ws.Range("P3").Value = 0.688622754491018
$ws.Range("S3").Value = 0.2095808383233533
$ws.Range("J4").Value = 0.02631578947368421
$ws.Range("P4").Value = 0.631578947368421
$ws.Range("S4").Value = 0.3421052631578947
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("S5").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.04694835680751173
$ws.Range("D6").Value = 0.009389671361502348
$ws.Range("F6").Value = 0.07511737089201878
$ws.Range("J6").Value = 0.2394366197183098
$ws.Range("O6").Value = 0.02347417840375587
$ws.Range("Q6").Value = 0.1549295774647887
$ws.Range("R6").Value = 0.05164319248826291
$ws.Range("S6").Value = 0.3990610328638498
$ws.Range("B7").Value = 0.09166666666666666
$ws.Range("D7").Value = 0.008333333333333333
$ws.Range("F7").Value = 0.05833333333333333
$ws.Range("J7").Value = 0.1375
$ws.Range("O7").Value = 0.008333333333333333
$ws.Range("Q7").Value = 0.175
$ws.Range("R7").Value = 0.075
$ws.Range("S7").Value = 0.4458333333333334
$ws.Range("B8").Value = 0.1147540983606557
$ws.Range("D8").Value = 0.01873536299765808
$ws.Range("F8").Value = 0.07728337236533958
$ws.Range("J8").Value = 0.09836065573770492
$ws.Range("O8").Value = 0.00702576112412178
$ws.Range("Q8").Value = 0.1826697892271663
$ws.Range("R8").Value = 0.06791569086651054
$ws.Range("S8").Value = 0.4332552693208431
$ws.Range("B9").Value = 0.1058823529411765
$ws.Range("D9").Value = 0.01176470588235294
$ws.Range("F9").Value = 0.07058823529411765
$ws.Range("J9").Value = 0.1117647058823529
$ws.Range("O9").Value = 0.005882352941176471
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.07058823529411765
$ws.Range("S9").Value = 0.4235294117647059
$ws.Range("B10").Value = 0.1119068934646374
$ws.Range("D10").Value = 0.02148612354521039
$ws.Range("E10").Value = 0.002685765443151298
$ws.Range("F10").Value = 0.07162041181736795
$ws.Range("J10").Value = 0.1298119964189794
$ws.Range("O10").Value = 0.01880035810205909
$ws.Range("Q10").Value = 0.1817367949865712
$ws.Range("R10").Value = 0.07162041181736795
$ws.Range("S10").Value = 0.3903312444046553
$ws.Range("G11").Value = 0.1471389645776567
$ws.Range("J11").Value = 0.05994550408719346
$ws.Range("K11").Value = 0.1825613079019074
$ws.Range("L11").Value = 0.5858310626702997
$ws.Range("S11").Value = 0.02452316076294278
$ws.Range("G12").Value = 0.7522935779816514
$ws.Range("J12").Value = 0.1467889908256881
$ws.Range("K12").Value = 0.004587155963302753
$ws.Range("L12").Value = 0.02752293577981652
$ws.Range("S12").Value = 0.06880733944954129
$ws.Range("G13").Value = 0.6818181818181818
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.06818181818181818
$ws.Range("J14").Value = 0.6
$ws.Range("S14").Value = 0.4
$ws.Range("F15").Value = 0.005128205128205128
$ws.Range("H15").Value = 0.1948717948717949
$ws.Range("I15").Value = 0.05128205128205128
$ws.Range("J15").Value = 0.3487179487179487
$ws.Range("K15").Value = 0.07179487179487179
$ws.Range("M15").Value = 0.01538461538461539
$ws.Range("N15").Value = 0.005128205128205128
$ws.Range("O15").Value = 0.04102564102564103
$ws.Range("S15").Value = 0.2666666666666667
$ws.Range("F16").Value = 0.01714285714285714
$ws.Range("H16").Value = 0.1485714285714286
$ws.Range("I16").Value = 0.06857142857142857
$ws.Range("J16").Value = 0.3714285714285714
$ws.Range("K16").Value = 0.1314285714285714
$ws.Range("M16").Value = 0.02857142857142857
$ws.Range("N16").Value = 0.005714285714285714
$ws.Range("O16").Value = 0.06857142857142857
$ws.Range("S16").Value = 0.16
$ws.Range("F17").Value = 0.01799485861182519
$ws.Range("H17").Value = 0.141388174807198
$ws.Range("I17").Value = 0.08483290488431877
$ws.Range("J17").Value = 0.4241645244215938
$ws.Range("K17").Value = 0.1182519280205656
$ws.Range("M17").Value = 0.01542416452442159
$ws.Range("O17").Value = 0.05912596401028278
$ws.Range("S17").Value = 0.1388174807197944
$ws.Range("F18").Value = 0.02040816326530612
$ws.Range("H18").Value = 0.1564625850340136
$ws.Range("I18").Value = 0.1156462585034014
$ws.Range("J18").Value = 0.3537414965986395
$ws.Range("K18").Value = 0.1292517006802721
$ws.Range("M18").Value = 0.02040816326530612
$ws.Range("O18").Value = 0.05442176870748299
$ws.Range("S18").Value = 0.1496598639455782
$ws.Range("F19").Value = 0.01780185758513932
$ws.Range("H19").Value = 0.2260061919504644
$ws.Range("I19").Value = 0.07739938080495357
$ws.Range("J19").Value = 0.3235294117647059
$ws.Range("K19").Value = 0.1462848297213622
$ws.Range("M19").Value = 0.021671826625387
$ws.Range("N19").Value = 0.002321981424148607
$ws.Range("O19").Value = 0.06578947368421052
$ws.Range("S19").Value = 0.1191950464396285
